$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the shared string text used in A18 ("Retest Application" -> "Regression Test")
$ws.Range("A18").Value = "Regression Test"

# Update the active cell selection on the sheet
$ws.Range("F5").Select()

# Update dates in column D for rows 16-19
$ws.Range("D16").Value = 44648
$ws.Range("D17").Value = 44648
$ws.Range("D18").Value = 44648
$ws.Range("D19").Value = 44641
